$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# --- 1. Capture the values currently in row 2 (before the insert shifts them down) ---
$oldB2 = $ws.Range("B2").Value2
$oldC2 = $ws.Range("C2").Value2
$oldD2 = $ws.Range("D2").Value2
$oldE2 = $ws.Range("E2").Value2
$oldF2 = $ws.Range("F2").Value2

# --- 2. Insert a brand-new row above the current row 2; everything below shifts down by one ---
$ws.Rows("2:2").Insert()

# Copy formatting (styles/number formats) from the (now shifted) data row 3 into the new row 2
$ws.Range("A3:F3").Copy()
$ws.Range("A2:F2").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# --- 3. Populate the new row 2 with the latest price entry ---
$ws.Range("A2").Value2 = "25-12-2025"
$ws.Range("B2").Value2 = $oldB2
$ws.Range("C2").Value2 = $oldC2
$ws.Range("D2").Value2 = $oldD2
$ws.Range("E2").Value2 = $oldE2
$ws.Range("F2").Value2 = $oldF2

# --- 4. Rebuild every hyperlink in the sheet so each F-cell's link target matches its text ---
$ws.Hyperlinks.Delete()

$lastRow = $ws.Range("A1048576").End(-4162).Row
for ($r = 2; $r -le $lastRow; $r++) {
    $cell = $ws.Cells.Item($r, 6)
    $text = $cell.Value2
    if ($text -and $text.ToString().StartsWith("http")) {
        $ws.Hyperlinks.Add($cell, $text.ToString())
    }
}
